$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "defaults"
$ws.Range("D2").Value = 1
$ws.Range("D3").Formula = "=1+D2"
$ws.Range("D4").Formula = "=1+D3"

$ws.Range("C2").Formula = "=B2/B3+D2"
$ws.Range("C3").Formula = "=C2*A2+D3"
$ws.Range("C4").Formula = "=B3^C2+D4"

$ws.Range("E5").Select()
